$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "stream:datastream"
$ws.Range("B2").Value = "dict"

$ws.Range("A3").Value = "operation_end_time"
$ws.Range("B3").Value = "datetime"

$ws.Range("A4").Value = "concept:name"
$ws.Range("B4").Value = "str"

$ws.Range("A5").Value = "SubProcessID"
$ws.Range("B5").Value = "str"

$ws.Range("A6").Value = "time:timestamp"
$ws.Range("B6").Value = "datetime"
